# Continue to GT data: append two labeled data blocks (rows 36-39) to Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 36: label "average" ---
$ws.Range("A36").Value = "average"

# --- Row 37: data row (uk/uk/uk placeholders in E:G) ---
$ws.Range("A37").Value = 3
$ws.Range("B37").Value = "GS"
$ws.Range("C37").Value = "rf all"
$ws.Range("D37").Value = 734

# --- Row 38: label "actual" ---
$ws.Range("A38").Value = "actual"

# --- Row 39: data row ---
$ws.Range("A39").Value = 3
$ws.Range("B39").Value = "GS"
# Populate the new unique string "All data test_20190401095702" before "uk"
# so shared-string table indices match the authored workbook ordering.
$ws.Range("C39").Value = "All data test_20190401095702"
$ws.Range("D39").Value = 734
$ws.Range("E39").Value = 333
$ws.Range("F39").Value = 401
$ws.Range("G39").Value = 33
$ws.Range("H39").Value = 0.90983606557377095
$ws.Range("I39").Value = 0.54632152588555905
$ws.Range("J39").Value = 0.83042394014962595
$ws.Range("K39").Value = 0.82763436753461594
$ws.Range("M39").Value = 2268

# Now fill the "uk" placeholders in row 37 (new unique string, will get next index)
$ws.Range("E37").Value = "uk"
$ws.Range("F37").Value = "uk"
$ws.Range("G37").Value = "uk"
$ws.Range("H37").Value = 0.91640722291407195
$ws.Range("I37").Value = 0.39818033122060897
$ws.Range("J37").Value = 1.5114249037227201
$ws.Range("K37").Value = 0.948787406201965
$ws.Range("M37").Value = 2268

# --- Update view state: scroll position and active selection ---
$ws.Range("E33").Select()
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1
